$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need a leading quote-prefix
# so Excel stores them as text (matching the inlineStr/text cells in the workbook),
# then ClearFormats() removes the quote-prefix style so the cell keeps the default style.

$ws.Range('D2').Value = '58.771.39'
$ws.Range('E2').Value = '  +0.26%  '

$ws.Range('D3').Value = '2.570.40'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = "'" + '560.83'
$ws.Range('E5').Value = '  +3.09%  '

$ws.Range('D6').Value = "'" + '142.24'
$ws.Range('E6').Value = '  -1.41%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = "'" + '0.591'
$ws.Range('E8').Value = '  +1.84%  '

$ws.Range('D9').Value = '2.575.89'
$ws.Range('E9').Value = '  -0.48%  '

$ws.Range('E10').Value = '  -1.91%  '

$ws.Range('E11').Value = '  +2.07%  '

$ws.Range('E12').Value = '  +8.36%  '

$ws.Range('D13').Value = "'" + '0.339'
$ws.Range('E13').Value = '  +1.98%  '

$ws.Range('D14').Value = '3.022.74'
$ws.Range('E14').Value = '  -0.47%  '

$ws.Range('D15').Value = '58.887.01'
$ws.Range('E15').Value = '  +0.60%  '

$ws.Range('D16').Value = "'" + '21.79'
$ws.Range('E16').Value = '  +5.55%  '

$ws.Range('E17').Value = '  +3.35%  '

$ws.Range('D18').Value = '2.573.21'
$ws.Range('E18').Value = '  -0.90%  '

$ws.Range('E19').Value = '  +0.77%  '

$ws.Range('D20').Value = "'" + '334.34'
$ws.Range('E20').Value = '  +0.01%  '

$ws.Range('D21').Value = "'" + '10.12'
$ws.Range('E21').Value = '  +0.59%  '

$ws.Range('E22').Value = '  +0.84%  '

$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('D24').Value = "'" + '63.92'
$ws.Range('E24').Value = '  -3.77%  '

$ws.Range('E25').Value = '  +4.45%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('E27').Value = '  +1.38%  '

$ws.Range('D28').Value = "'" + '7.18'
$ws.Range('E28').Value = '  +1.02%  '

$ws.Range('D29').Value = '0.0₃0776'
$ws.Range('E29').Value = '  +4.74%  '

$ws.Range('D30').Value = "'" + '0.999'
$ws.Range('E30').Value = '  +0.03%  '

$ws.Range('E31').Value = '  +1.52%  '

$ws.Range('D32').Value = "'" + '160.45'
$ws.Range('E32').Value = '  +5.02%  '

$ws.Range('D33').Value = "'" + '6.01'
$ws.Range('E33').Value = '  +0.89%  '

$ws.Range('D34').Value = "'" + '18.84'
$ws.Range('E34').Value = '  -0.36%  '

$ws.Range('D35').Value = "'" + '3.98'
$ws.Range('E35').Value = '  +1.89%  '

$ws.Range('D36').Value = "'" + '0.874'
$ws.Range('E36').Value = '  +2.66%  '

$ws.Range('D37').Value = "'" + '0.874'
$ws.Range('E37').Value = '  +6.34%  '

$ws.Range('D38').Value = "'" + '1.12'
$ws.Range('E38').Value = '  +2.38%  '

$ws.Range('D39').Value = "'" + '36.72'
$ws.Range('E39').Value = '  -1.19%  '

$ws.Range('E40').Value = '  +3.46%  '

$ws.Range('D41').Value = "'" + '293.09'
$ws.Range('E41').Value = '  +4.97%  '

$ws.Range('D42').Value = "'" + '3.59'
$ws.Range('E42').Value = '  +0.30%  '

$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('D44').Value = "'" + '0.0968'
$ws.Range('E44').Value = '  +2.86%  '

$ws.Range('D45').Value = "'" + '0.592'
$ws.Range('E45').Value = '  -0.16%  '

$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = "'" + '10.61'
$ws.Range('E46').Value = '  -0.19%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = "'" + '0.0533'
$ws.Range('E47').Value = '  +0.99%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'" + '18.91'
$ws.Range('E48').Value = '  +1.97%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = "'" + '124.59'
$ws.Range('E49').Value = '  +13.69%  '

$ws.Range('D50').Value = "'" + '0.0230'
$ws.Range('E50').Value = '  +1.15%  '

$ws.Range('D51').Value = "'" + '18.35'
$ws.Range('E51').Value = '  +2.48%  '

# Remove the quote-prefix number format so styling matches the original plain text cells
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D33').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
